# Updated symbol list on Wed Jan  4 05:23:21 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns on the crypto symbol sheet.
# Values are stored as plain text in the source data (e.g. "255.24", "3.96%"),
# so we force each target cell to Text format before writing the new value.
# This preserves exact formatting (trailing zeros, percent signs, thousands
# separators, etc.) instead of letting Excel silently re-interpret the text
# as a number/percentage and round/reformat it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "255.24";    "E2"  = "3.96%";
    "D3"  = "28.20";     "E3"  = "-3.99%";
                         "E4"  = "4.55%";
    "D5"  = "0.05827";   "E5"  = "0.73%";
    "D6"  = "6.722";     "E6"  = "1.61%";
    "D7"  = "0.8667";    "E7"  = "1.59%";
    "D8"  = "0.9091";    "E8"  = "5.61%";
    "D9"  = "0.1418";    "E9"  = "4.17%";
    "D10" = "0.07179";   "E10" = "1.69%";
    "D11" = "0.03202";   "E11" = "1.09%";
                         "E12" = "-1.29%";
                         "E13" = "0.21%";
    "D14" = "0.0006040"; "E14" = "-94.11%";
    "D15" = "0.005939";  "E15" = "-1.77%";
    "D16" = "3.495";     "E16" = "0.26%";
    "D18" = "2.273";     "E18" = "5.20%";
    "D19" = "0.3169";    "E19" = "-0.97%";
                         "E20" = "4.42%";
                         "E21" = "2.24%";
    "D22" = "3.529";     "E22" = "6.55%";
    "D23" = "0.04160";   "E23" = "0.56%";
                         "E24" = "-1.40%";
    "D25" = "0.001228";  "E25" = "0.21%";
    "D26" = "0.004862";  "E26" = "17.59%";
                         "E27" = "-0.75%";
    "D28" = "0.0001938"; "E28" = "34.23%";
    "D40" = "0.03844";   "E40" = "3.05%";
    "D41" = "0.005739";  "E41" = "-0.45%";
    "D42" = "0.1098";    "E42" = "2.82%";
                         "E43" = "0.08%";
    "D44" = "0.009852";  "E44" = "7.43%";
    "D45" = "0.00005282";"E45" = "0.16%";
                         "E46" = "0.08%";
    "D47" = "0.1001";    "E47" = "72.75%";
    "D48" = "0.002213";  "E48" = "1.85%";
                         "E49" = "0.08%";
                         "E50" = "0.08%";
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
